# Generate Report for Handback
# Update the "last generated / handback" timestamps and the zh-cn/de-de
# "ht" -> "mt" status cells.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G)
$overview.Range("G2").Value = "2016-08-30 04:15:58"
$overview.Range("G4").Value = "2016-08-30 04:15:58"

# zh-cn sheet
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E4").Value = "mt"
$zhcn.Range("H2").Value = "2016-08-30 04:15:53"
$zhcn.Range("H4").Value = "2016-08-30 04:15:53"
$zhcn.Range("K2").Value = "2016-08-30 04:16:14"
$zhcn.Range("K4").Value = "2016-08-30 04:16:14"

# de-de sheet
$dede.Range("E2").Value = "mt"
$dede.Range("E4").Value = "mt"
$dede.Range("H2").Value = "2016-08-30 04:15:58"
$dede.Range("H4").Value = "2016-08-30 04:15:58"
$dede.Range("K2").Value = "2016-08-30 04:16:21"
$dede.Range("K4").Value = "2016-08-30 04:16:21"
